# edit.ps1 - apply the recorded edit to the open presentation
#
# The source commit:
#  1) changes the table style GUID used by the table on slide 16
#     (shape 3, the graphicFrame holding the table) from
#     {681CF0FF-52AA-4FEF-862C-B0B93F731B36} to
#     {CD90D8F6-83DD-48A7-B402-B2A74260EF36}.
#  2) swaps the deck's theme colors: the custom "Integral" palette is
#     replaced by the stock "Office" palette (what had been sitting,
#     unused, on the Notes Master's theme).

$p = $ppt.ActivePresentation

# --- 1) table style -------------------------------------------------
$s   = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table
$tbl.ApplyStyle("{CD90D8F6-83DD-48A7-B402-B2A74260EF36}")

# --- 2) theme colors --------------------------------------------------
# Order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink. Values are standard VBA RGB() integers
# (R + G*256 + B*65536) for the stock Office theme palette.
$officeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeRGB[$i - 1]
}
